# Clases.xlsx — "lucia cambia xls de clases"
#
# 1) Fix the Gerber & Green / Barabas reading-list cell (C2): the line break
#    between the two references is replaced by a single space, and the stray
#    space right after "[Descarga]" is removed.
# 2) Move the active selection from C2 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0
$newText = "Gerber " + $nbsp + "& Green 2012. FEDAI [Descarga]https://drive.google.com/drive/folders/14HDw0lx7v8cduNtj2XNvvZ5fm_lQ7Z6y?usp=sharing) Barabas 2010 [pdf](https://drive.google.com/u/0/uc?id=15SqCaheQIA_Eg8Q6CxkkF5Gdt2dPdK1Y&export=download)  Aronow et al 2015 [pdf]()"

$ws.Range("C2").Value = $newText

$ws.Range("C3").Select()
